# Insert a new data row at row 97 (pushing the existing rows 97-168 down to 98-169)
# and populate it with a new price-report entry for "Jengibre" (Vega Modelo de Temuco).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 97; this shifts rows 97..168 down to 98..169
# and Excel automatically grows the sheet dimension to A1:R169.
$ws.Rows.Item(97).Insert()

# Fill in the new row 97 with the reported values.
$ws.Cells.Item(97, 1).Value  = 10
$ws.Cells.Item(97, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(97, 3).Value  = "La Araucanía"
$ws.Cells.Item(97, 4).Value  = 44719
$ws.Cells.Item(97, 5).Value  = 9
$ws.Cells.Item(97, 6).Value  = 100114007
$ws.Cells.Item(97, 7).Value  = "Jengibre"
$ws.Cells.Item(97, 8).Value  = "Sin especificar"
$ws.Cells.Item(97, 9).Value  = "Primera"
$ws.Cells.Item(97, 10).Value = 50
$ws.Cells.Item(97, 11).Value = 20000
$ws.Cells.Item(97, 12).Value = 20000
$ws.Cells.Item(97, 13).Value = 20000
$ws.Cells.Item(97, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(97, 15).Value = "Perú"
$ws.Cells.Item(97, 16).Value = 1538
$ws.Cells.Item(97, 17).Value = 13
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# Match the date-number format used by the other rows in column D.
$ws.Cells.Item(97, 4).NumberFormat = $ws.Cells.Item(98, 4).NumberFormat
